$wb = $excel.ActiveWorkbook

# Work on the "RT" worksheet (sheet2) and type headers into row 1.
$rt = $wb.Worksheets.Item("RT")

# Typing order observed from the resulting shared-string table:
# A1, B1, C1, E1, D1 (column D typed last, then selected as a whole column).
$rt.Range("A1").Value = "Trial"
$rt.Range("B1").Value = "Abort (2651)"
$rt.Range("C1").Value = "CorrectSaccade (2600)"
$rt.Range("E1").Value = "RT"
$rt.Range("D1").Value = "incorrectSaccade (887)"

# Column widths, matching Excel's auto "best fit" for the typed headers.
$rt.Columns.Item(1).ColumnWidth = 7.5
$rt.Columns.Item(2).ColumnWidth = 10.666666666666668
$rt.Columns.Item(3).ColumnWidth = 18.5
$rt.Columns.Item(4).ColumnWidth = 14.166666666666668
$rt.Columns.Item(5).ColumnWidth = 5.666666666666666

# Make RT the active (and selected) sheet/tab.
$rt.Activate()

# Select the whole column D as the last interactive action on this sheet.
$rt.Range("D1:D1048576").Select() | Out-Null
